$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record was added to the daily series. In the source data it
# lands at row 248, pushing every existing row from 248 downward by one
# (old row 351 becomes new row 352).
$ws.Rows("248:248").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

$ws.Range("A248").Value = 7
$ws.Range("B248").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C248").Value = "Ñuble"
$ws.Range("D248").Value = 44755
$ws.Range("E248").Value = 16
$ws.Range("F248").Value = 100114001
$ws.Range("G248").Value = "Papa"
$ws.Range("H248").Value = "Asterix"
$ws.Range("I248").Value = "1a (guarda)"
$ws.Range("J248").Value = 200
$ws.Range("K248").Value = 7000
$ws.Range("L248").Value = 7500
$ws.Range("M248").Value = 7250
$ws.Range("N248").Value = "$/saco 25 kilos"
$ws.Range("O248").Value = "Provincia de Diguillín"
$ws.Range("P248").Value = 290
$ws.Range("Q248").Value = 25
$ws.Range("R248").Value = "Hortaliza"
